$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F12").Formula = "=E12-D12"
Write-Host "done"
